$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "36.724.85"
Set-TextValue $ws "E2" "  -0.86%  "

Set-TextValue $ws "D3" "2.060.89"
Set-TextValue $ws "E3" "  +0.21%  "

Set-TextValue $ws "E4" "  +0.02%  "

Set-TextValue $ws "D5" "246.71"
Set-TextValue $ws "E5" "  +0.06%  "

Set-TextValue $ws "E6" "  +1.32%  "

Set-TextValue $ws "D8" "55.35"
Set-TextValue $ws "E8" "  -6.62%  "

Set-TextValue $ws "D9" "60.93"
Set-TextValue $ws "E9" "  +2.56%  "

Set-TextValue $ws "D10" "0.369"
Set-TextValue $ws "E10" "  -2.88%  "

Set-TextValue $ws "D11" "0.0754"
Set-TextValue $ws "E11" "  -2.81%  "

Set-TextValue $ws "E12" "  -2.94%  "

Set-TextValue $ws "D13" "0.968"
Set-TextValue $ws "E13" "  +8.79%  "

Set-TextValue $ws "D14" "14.87"
Set-TextValue $ws "E14" "  -4.36%  "

Set-TextValue $ws "D15" "2.362.75"
Set-TextValue $ws "E15" "  +0.18%  "

Set-TextValue $ws "D16" "5.49"
Set-TextValue $ws "E16" "  -4.40%  "

Set-TextValue $ws "D17" "2.054.98"
Set-TextValue $ws "E17" "  -0.71%  "

Set-TextValue $ws "D18" "36.671.18"
Set-TextValue $ws "E18" "  -0.95%  "

Set-TextValue $ws "D19" "17.34"
Set-TextValue $ws "E19" "  -5.59%  "

Set-TextValue $ws "D20" "72.47"
Set-TextValue $ws "E20" "  -2.40%  "

Set-TextValue $ws "E21" "  -3.31%  "

Set-TextValue $ws "D22" "238.96"
Set-TextValue $ws "E22" "  -0.09%  "

Set-TextValue $ws "D23" "5.28"
Set-TextValue $ws "E23" "  -3.55%  "

Set-TextValue $ws "E24" "  -0.05%  "

Set-TextValue $ws "E25" "  -2.47%  "

Set-TextValue $ws "D26" "2.27"
Set-TextValue $ws "E26" "  +5.05%  "

Set-TextValue $ws "D27" "9.29"
Set-TextValue $ws "E27" "  -6.74%  "

Set-TextValue $ws "D28" "166.12"
Set-TextValue $ws "E28" "  -2.11%  "

Set-TextValue $ws "E29" "  +0.17%  "

Set-TextValue $ws "E30" "  -1.32%  "

Set-TextValue $ws "D31" "1.21"
Set-TextValue $ws "E31" "  +8.55%  "

Set-TextValue $ws "D32" "5.10"
Set-TextValue $ws "E32" "  -6.73%  "

Set-TextValue $ws "D33" "4.54"
Set-TextValue $ws "E33" "  -3.68%  "

Set-TextValue $ws "E34" "  -3.39%  "

Set-TextValue $ws "E35" "  +0.05%  "

Set-TextValue $ws "B36" "LidoDAOToken"
Set-TextValue $ws "C36" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws "D36" "2.28"
Set-TextValue $ws "E36" "  -2.83%  "

Set-TextValue $ws "B37" "WEMIXToken"
Set-TextValue $ws "C37" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws "D37" "1.84"
Set-TextValue $ws "E37" "  +0.61%  "

Set-TextValue $ws "B38" "Kaspa"
Set-TextValue $ws "C38" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D38" "0.0854"
Set-TextValue $ws "E38" "  +1.64%  "

Set-TextValue $ws "E39" "  -3.13%  "

Set-TextValue $ws "E40" "  -5.63%  "

Set-TextValue $ws "E41" "  -5.66%  "

Set-TextValue $ws "D42" "0.0217"
Set-TextValue $ws "E42" "  -3.75%  "

Set-TextValue $ws "E43" "  -4.82%  "

Set-TextValue $ws "D44" "95.32"
Set-TextValue $ws "E44" "  -2.79%  "

Set-TextValue $ws "E45" "  -4.25%  "

Set-TextValue $ws "D46" "1.422.25"
Set-TextValue $ws "E46" "  +8.93%  "

Set-TextValue $ws "D47" "7.67"
Set-TextValue $ws "E47" "  +12.10%  "

Set-TextValue $ws "E48" "  -6.27%  "

Set-TextValue $ws "D49" "2.94"
Set-TextValue $ws "E49" "  +2.52%  "

Set-TextValue $ws "D50" "2.29"
Set-TextValue $ws "E50" "  -4.08%  "

Set-TextValue $ws "D51" "2.248.28"
Set-TextValue $ws "E51" "  +0.00%  "
